$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48. This shifts existing rows 48..113 down to 49..114
# and extends the used range / dimension accordingly.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new record's data.
$ws.Cells.Item(48, 1).Value = 3
$ws.Cells.Item(48, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(48, 3).Value = "Coquimbo"
$ws.Cells.Item(48, 4).Value = (Get-Date -Year 2021 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(48, 5).Value = 5
$ws.Cells.Item(48, 6).Value = 100112001
$ws.Cells.Item(48, 7).Value = "Berenjena"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 105
$ws.Cells.Item(48, 11).Value = 11500
$ws.Cells.Item(48, 12).Value = 12000
$ws.Cells.Item(48, 13).Value = 11762
$ws.Cells.Item(48, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 196
$ws.Cells.Item(48, 17).Value = 60
$ws.Cells.Item(48, 18).Value = "Hortaliza"
